$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 border cell
$ws.Range("J3").Value = ""

# Row 4 header year
$ws.Range("J4").Value = 2021

# Row 5-6 (customFormat rows)
$ws.Range("J5").Value = 5356.3
$ws.Range("J6").Value = 9.5

# Rows 8-17 (plain data rows, style picks up column default)
$ws.Range("J8").Value = 7.9
$ws.Range("J9").Value = 10.5
$ws.Range("J11").Value = 9.6
$ws.Range("J12").Value = 9.4
$ws.Range("J14").Value = 14.8
$ws.Range("J15").Value = 9.1
$ws.Range("J16").Value = 9.5
$ws.Range("J17").Value = 5.9

# Rows 19-26 need numFmt 164 (0.0) style group
$ws.Range("J19").Value = 12.434613462352335
$ws.Range("J20").Value = 16.80050595536094
$ws.Range("J21").Value = 11.282963378125267
$ws.Range("J22").Value = 25.042808754677555
$ws.Range("J23").Value = 3.2011163356916352
$ws.Range("J24").Value = 13.523574517571838
$ws.Range("J25").Value = 6.1196997869329204
$ws.Range("J26").Value = 5.9488136666578013

# Row 27 (bottom border row)
$ws.Range("J27").Value = 5.2451982064110645

$ws.Range("N8").Select()
